$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to be stored as a plain text
# value (matching the workbook's existing inline-string/shared-string cells)
# instead of letting Excel auto-coerce numeric-looking text into a number.
# The cell's NumberFormat/Style is restored to the default ("Normal") right
# after the write so no visible/structural formatting change is introduced.
function Set-TextValue($sheet, $addr, $val) {
    $c = $sheet.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue $ws "D2"  "245.70"
Set-TextValue $ws "D3"  "23.98"
Set-TextValue $ws "D4"  "5.351"
Set-TextValue $ws "D5"  "0.05807"
Set-TextValue $ws "D6"  "3.372"
Set-TextValue $ws "D7"  "6.473"
Set-TextValue $ws "D9"  "0.9200"
Set-TextValue $ws "D10" "0.1404"
Set-TextValue $ws "D11" "0.07398"
Set-TextValue $ws "D12" "0.03190"
Set-TextValue $ws "D13" "0.03066"
Set-TextValue $ws "D14" "0.09381"
Set-TextValue $ws "D15" "3.849"
Set-TextValue $ws "D16" "0.001577"
Set-TextValue $ws "D17" "0.04696"
Set-TextValue $ws "D18" "0.0005967"

# --- Row 18 (One / ONE) also gets a "Worstin24h" suffix on its summary cell ---
Set-TextValue $ws "E18" "17OneONEWorstin24h"

Set-TextValue $ws "D19" "0.005915"
Set-TextValue $ws "D20" "0.001244"
Set-TextValue $ws "D22" "0.00008794"
Set-TextValue $ws "D23" "3.596"
Set-TextValue $ws "D25" "0.3184"
Set-TextValue $ws "D26" "0.1320"
Set-TextValue $ws "D28" "0.0002349"
Set-TextValue $ws "D40" "0.03841"

# --- Rows 41-43: the coin ranking rotated (Kick/BKEX/CEJI shifted down a slot) ---
Set-TextValue $ws "B41" "KickToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006326"
Set-TextValue $ws "E41" "40KickTokenKICK"

Set-TextValue $ws "B42" "BKEXToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1066"
Set-TextValue $ws "E42" "41BKEXTokenBKK"

Set-TextValue $ws "B43" "CEJI"
Set-TextValue $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002598"
Set-TextValue $ws "E43" "42CEJICEJI"

Set-TextValue $ws "D44" "0.009044"
Set-TextValue $ws "D45" "0.00005248"
Set-TextValue $ws "D47" "0.6852"
Set-TextValue $ws "D48" "0.001831"
Set-TextValue $ws "D49" "0.00002099"
Set-TextValue $ws "D50" "0.0001999"
